$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Populate cell values.
#    The order below matters: it reproduces the exact order in which new
#    shared-string entries were first introduced (Source_URL, Target_URL,
#    the target URL value, hemag, repo, tests, Target_Pat) while reusing the
#    already-present strings (Source_Username, Source_Pat, Target_Username,
#    http://172.191.4.85/DefaultCollection, lexcon, xs6md..., BqV9Eb...).
# ---------------------------------------------------------------------------

# Row 1 headers (A:E first)
$ws.Range("A1").Value = "Source_URL"
$ws.Range("B1").Value = "Source_Username"
$ws.Range("C1").Value = "Source_Pat"
$ws.Range("D1").Value = "Target_URL"
$ws.Range("E1").Value = "Target_Username"

# Row 2 data (A:E)
$ws.Range("A2").Value = "http://172.191.4.85/DefaultCollection"
$ws.Range("B2").Value = "lexcon"
$ws.Range("C2").Value = "xs6mdazt46rfd2iur7nrpmkyprnsnvpxyizehr5yltliv3aaacaa"
$ws.Range("D2").Value = "https://dev.azure.com/PLMigration"
$ws.Range("E2").Value = "hemag"

# New repo/tests columns
$ws.Range("G1").Value = "repo"
$ws.Range("G2").Value = "tests"

# Target_Pat header + secret value added last
$ws.Range("F1").Value = "Target_Pat"
$ws.Range("F2").Value = "BqV9EbVuxxzXtmzEdtdTfevv1qZ3EQszfR410EtLL0TDvbwxMruhJQQJ99AKACAAAAAyb0Q7AAASAZDOBwc3"

# ---------------------------------------------------------------------------
# 2. Font formatting - explicit black font color on the regular data cells.
# ---------------------------------------------------------------------------

# Column E previously carried wrapText - the new layout no longer wraps it,
# so reset to the Normal style before re-applying the font color.
$ws.Range("E1:E2").Style = "Normal"

$ws.Range("A1").Font.Color = 0
$ws.Range("B1").Font.Color = 0
$ws.Range("C1").Font.Color = 0
$ws.Range("D1").Font.Color = 0
$ws.Range("E1").Font.Color = 0
$ws.Range("G1").Font.Color = 0

$ws.Range("B2").Font.Color = 0
$ws.Range("C2").Font.Color = 0
$ws.Range("E2").Font.Color = 0
$ws.Range("G2").Font.Color = 0

# Target_Pat column (F) - black font plus wrap text.
$ws.Range("F1:F2").Font.Color = 0
$ws.Range("F1:F2").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Hyperlinks.
#    A2 already carried a hyperlink to http://172.191.4.85/ - recreate it so
#    it also gets a display tooltip, then restore the cell's full text.
#    D2 is a brand-new hyperlink pointing at the Azure DevOps org URL.
# ---------------------------------------------------------------------------
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "http://172.191.4.85/", [Type]::Missing, [Type]::Missing, "http://172.191.4.85/") | Out-Null
$ws.Range("A2").Value = "http://172.191.4.85/DefaultCollection"

$ws.Hyperlinks.Add($ws.Range("D2"), "https://dev.azure.com/PLMigration") | Out-Null

# ---------------------------------------------------------------------------
# 4. Column widths for the two newly introduced columns.
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 34.5
$ws.Columns("F").ColumnWidth = 99

# ---------------------------------------------------------------------------
# 5. Selection.
# ---------------------------------------------------------------------------
$ws.Range("B2").Select() | Out-Null
